# Update TPM-derived NATMI ligand-receptor values on Sheet1
# Sets the recomputed values for columns G-J (ligand stats), M-P (receptor stats),
# and Q-T (edge weights/specificities) across data rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.458056666666667
$ws.Range("H2").Value = 4.374169999999999
$ws.Range("I2").Value = 0.2323568509805328
$ws.Range("J2").Value = 0.2323568509805327
$ws.Range("M2").Value = 8.141899
$ws.Range("N2").Value = 24.425697
$ws.Range("O2").Value = 0.1319208574566759
$ws.Range("P2").Value = 0.1319208574566759
$ws.Range("Q2").Value = 11.87135011627667
$ws.Range("R2").Value = 106.84215104649
$ws.Range("S2").Value = 0.03065271501728495
$ws.Range("T2").Value = 0.03065271501728495
$ws.Range("G3").Value = 1.458056666666667
$ws.Range("H3").Value = 4.374169999999999
$ws.Range("I3").Value = 0.2323568509805328
$ws.Range("J3").Value = 0.2323568509805327
$ws.Range("O3").Value = 0.2680499808994311
$ws.Range("P3").Value = 0.2680499808994311
$ws.Range("Q3").Value = 24.12139545836
$ws.Range("R3").Value = 217.09255912524
$ws.Range("S3").Value = 0.06228324946718377
$ws.Range("T3").Value = 0.06228324946718376
$ws.Range("G4").Value = 1.458056666666667
$ws.Range("H4").Value = 4.374169999999999
$ws.Range("I4").Value = 0.2323568509805328
$ws.Range("J4").Value = 0.2323568509805327
$ws.Range("M4").Value = 11.387163
$ws.Range("N4").Value = 34.161489
$ws.Range("O4").Value = 0.1845029405251692
$ws.Range("P4").Value = 0.1845029405251692
$ws.Range("Q4").Value = 16.60312892657
$ws.Range("R4").Value = 149.42816033913
$ws.Range("S4").Value = 0.04287052225707683
$ws.Range("T4").Value = 0.04287052225707683
$ws.Range("G5").Value = 1.458056666666667
$ws.Range("H5").Value = 4.374169999999999
$ws.Range("I5").Value = 0.2323568509805328
$ws.Range("J5").Value = 0.2323568509805327
$ws.Range("M5").Value = 10.07930033333333
$ws.Range("N5").Value = 30.237901
$ws.Range("O5").Value = 0.163312016341236
$ws.Range("P5").Value = 0.163312016341236
$ws.Range("Q5").Value = 14.69619104635222
$ws.Range("R5").Value = 132.26571941717
$ws.Range("S5").Value = 0.0379466658443309
$ws.Range("T5").Value = 0.0379466658443309
$ws.Range("G6").Value = 1.458056666666667
$ws.Range("H6").Value = 4.374169999999999
$ws.Range("I6").Value = 0.2323568509805328
$ws.Range("J6").Value = 0.2323568509805327
$ws.Range("M6").Value = 15.56617066666666
$ws.Range("N6").Value = 46.69851199999999
$ws.Range("O6").Value = 0.2522142047774878
$ws.Range("P6").Value = 0.2522142047774878
$ws.Range("Q6").Value = 22.69635891500444
$ws.Range("R6").Value = 204.26723023504
$ws.Range("S6").Value = 0.05860369839465631
$ws.Range("T6").Value = 0.0586036983946563
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3032511962008422
$ws.Range("J7").Value = 0.3032511962008422
$ws.Range("M7").Value = 8.141899
$ws.Range("N7").Value = 24.425697
$ws.Range("O7").Value = 0.1319208574566759
$ws.Range("P7").Value = 0.1319208574566759
$ws.Range("Q7").Value = 15.493415012676
$ws.Range("R7").Value = 139.440735114084
$ws.Range("S7").Value = 0.04000515782757777
$ws.Range("T7").Value = 0.04000515782757777
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3032511962008422
$ws.Range("J8").Value = 0.3032511962008422
$ws.Range("O8").Value = 0.2680499808994311
$ws.Range("P8").Value = 0.2680499808994311
$ws.Range("Q8").Value = 31.481068864176
$ws.Range("S8").Value = 0.08128647734936539
$ws.Range("T8").Value = 0.08128647734936539
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3032511962008422
$ws.Range("J9").Value = 0.3032511962008422
$ws.Range("M9").Value = 11.387163
$ws.Range("N9").Value = 34.161489
$ws.Range("O9").Value = 0.1845029405251692
$ws.Range("P9").Value = 0.1845029405251692
$ws.Range("Q9").Value = 21.668905764612
$ws.Range("R9").Value = 195.020151881508
$ws.Range("S9").Value = 0.0559507374168304
$ws.Range("T9").Value = 0.05595073741683041
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3032511962008422
$ws.Range("J10").Value = 0.3032511962008422
$ws.Range("M10").Value = 10.07930033333333
$ws.Range("N10").Value = 30.237901
$ws.Range("O10").Value = 0.163312016341236
$ws.Range("P10").Value = 0.163312016341236
$ws.Range("Q10").Value = 19.180142507508
$ws.Range("R10").Value = 172.621282567572
$ws.Range("S10").Value = 0.0495245643094513
$ws.Range("T10").Value = 0.0495245643094513
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3032511962008422
$ws.Range("J11").Value = 0.3032511962008422
$ws.Range("M11").Value = 15.56617066666666
$ws.Range("N11").Value = 46.69851199999999
$ws.Range("O11").Value = 0.2522142047774878
$ws.Range("P11").Value = 0.2522142047774878
$ws.Range("Q11").Value = 29.621239749696
$ws.Range("R11").Value = 266.5911577472639
$ws.Range("S11").Value = 0.07648425929761735
$ws.Range("T11").Value = 0.07648425929761735
$ws.Range("G12").Value = 2.914094333333333
$ws.Range("H12").Value = 8.742283
$ws.Range("I12").Value = 0.4643919528186251
$ws.Range("J12").Value = 0.4643919528186251
$ws.Range("M12").Value = 8.141899
$ws.Range("N12").Value = 24.425697
$ws.Range("O12").Value = 0.1319208574566759
$ws.Range("P12").Value = 0.1319208574566759
$ws.Range("Q12").Value = 23.72626173847233
$ws.Range("R12").Value = 213.536355646251
$ws.Range("S12").Value = 0.0612629846118132
$ws.Range("T12").Value = 0.0612629846118132
$ws.Range("G13").Value = 2.914094333333333
$ws.Range("H13").Value = 8.742283
$ws.Range("I13").Value = 0.4643919528186251
$ws.Range("J13").Value = 0.4643919528186251
$ws.Range("O13").Value = 0.2680499808994311
$ws.Range("P13").Value = 0.2680499808994311
$ws.Range("Q13").Value = 48.209389541764
$ws.Range("R13").Value = 433.884505875876
$ws.Range("S13").Value = 0.124480254082882
$ws.Range("T13").Value = 0.124480254082882
$ws.Range("G14").Value = 2.914094333333333
$ws.Range("H14").Value = 8.742283
$ws.Range("I14").Value = 0.4643919528186251
$ws.Range("J14").Value = 0.4643919528186251
$ws.Range("M14").Value = 11.387163
$ws.Range("N14").Value = 34.161489
$ws.Range("O14").Value = 0.1845029405251692
$ws.Range("P14").Value = 0.1845029405251692
$ws.Range("Q14").Value = 33.18326717104301
$ws.Range("R14").Value = 298.649404539387
$ws.Range("S14").Value = 0.08568168085126195
$ws.Range("T14").Value = 0.08568168085126196
$ws.Range("G15").Value = 2.914094333333333
$ws.Range("H15").Value = 8.742283
$ws.Range("I15").Value = 0.4643919528186251
$ws.Range("J15").Value = 0.4643919528186251
$ws.Range("M15").Value = 10.07930033333333
$ws.Range("N15").Value = 30.237901
$ws.Range("O15").Value = 0.163312016341236
$ws.Range("P15").Value = 0.163312016341236
$ws.Range("Q15").Value = 29.37203198533145
$ws.Range("R15").Value = 264.348287867983
$ws.Range("S15").Value = 0.07584078618745378
$ws.Range("T15").Value = 0.07584078618745378
$ws.Range("G16").Value = 2.914094333333333
$ws.Range("H16").Value = 8.742283
$ws.Range("I16").Value = 0.4643919528186251
$ws.Range("J16").Value = 0.4643919528186251
$ws.Range("M16").Value = 15.56617066666666
$ws.Range("N16").Value = 46.69851199999999
$ws.Range("O16").Value = 0.2522142047774878
$ws.Range("P16").Value = 0.2522142047774878
$ws.Range("Q16").Value = 45.36128973143288
$ws.Range("R16").Value = 408.251607582896
$ws.Range("S16").Value = 0.1171262470852141
$ws.Range("T16").Value = 0.1171262470852141

$wb.Save()
